$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("navngivning")

# Order matters for shared-string table ordering, mirror the order new
# strings were introduced: C20, F20, C19, B20, B19
$ws.Range("C20").Value = "FS-OC1.1 : seHistorik"

$ws.Range("F20").Value = "remember to check Anders ppt (OOA-6)"
$ws.Range("F20").Interior.Color = $ws.Range("B7").Interior.Color
$ws.Range("F20").HorizontalAlignment = $ws.Range("B7").HorizontalAlignment
$ws.Range("F20").VerticalAlignment = $ws.Range("B7").VerticalAlignment
$ws.Range("F20").WrapText = $ws.Range("B7").WrapText

$ws.Range("C19").Value = "FS-UC1 : seHistorik"
$ws.Range("B20").Value = "FS-OC#.# : navn"
$ws.Range("B19").Value = "FS-UC# : navn"

# Column F width change
$ws.Columns.Item(6).ColumnWidth = 16.44140625

# Selection change
$ws.Range("D20").Select()
